$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "_old" / "_new" header labels to "_FV2310" / "_FV2404" ---
# Columns A..J (1-10) carry the "_old" suffixed headers.
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}

# Column K (11) is "diff" - untouched.

# Columns L..U (12-21) carry the "_new" suffixed headers.
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# --- 2. Turn the used range into an Excel Table (ListObject) named "Table1" ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U84"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split below row 1, frozen) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
